# This script reproduces the commit "Code in extent branch with extent work",
# which records that the "Navigate" test scenario (rows 2-8 of Test Steps /
# row 2 of Test Cases) finished with result "Pass" (re-cased from "PASS"),
# and that the "Login" test scenario (rows 9-15 of Test Steps / row 3 of
# Test Cases) also finished running: RunMode flipped from "No" to "Yes" and
# its Results column got populated with "Pass" for every step.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Test Cases")
$ws2 = $wb.Worksheets.Item("Test Steps")

# ---------------------------------------------------------------------------
# 1. Cell content changes
# ---------------------------------------------------------------------------

# "Test Cases" sheet
#   D2 (Navigate scenario Results): "PASS" -> "Pass"
#   C3 (Login scenario RunMode):    "No"   -> "Yes"
#   D3 (Login scenario Results):    ""     -> "Pass"
$ws1.Range("D2").Value = "Pass"
$ws1.Range("C3").Value = "Yes"
$ws1.Range("D3").Value = "Pass"

# "Test Steps" sheet
#   G2:G8  (Navigate scenario step results): "PASS" -> "Pass"
#   G9:G15 (Login scenario step results):    ""     -> "Pass"
$ws2.Range("G2:G8").Value = "Pass"
$ws2.Range("G9:G15").Value = "Pass"

# ---------------------------------------------------------------------------
# 2. Column widths for the (now fully populated) "Results" column, and the
#    narrower "TestCase ID" column on the Test Steps sheet.
# ---------------------------------------------------------------------------

$ws1.Columns.Item(4).ColumnWidth = 8.75   # Results column on Test Cases
$ws2.Columns.Item(7).ColumnWidth = 8.75   # Results column on Test Steps
$ws2.Columns.Item(1).ColumnWidth = 13.6   # TestCase ID column on Test Steps

# ---------------------------------------------------------------------------
# 3. Active sheet / selection bookkeeping
#    The workbook now opens on "Test Steps" (activeTab = 1) with H22
#    selected, while "Test Cases" keeps a remembered selection of I14.
# ---------------------------------------------------------------------------

$ws1.Activate()
$ws1.Range("I14").Select()

$ws2.Activate()
$ws2.Range("H22").Select()

Write-Host "Applied extent-work result updates to Navigate workbook"
